$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 232, shifting the rest of the
# table (old rows 232-337) down by one (new rows 233-338).
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with the new price-quote record.
$ws.Range("A232").Value = 5
$ws.Range("B232").Value = "Macroferia Regional de Talca"
$ws.Range("C232").Value = "Maule"
$ws.Range("D232").Value = 44960
$ws.Range("E232").Value = 7
$ws.Range("F232").Value = "Fruta"
$ws.Range("G232").Value = 100108
$ws.Range("H232").Value = "Tropicales y subtropicales"
$ws.Range("I232").Value = 100108005
$ws.Range("J232").Value = "Piña"
$ws.Range("K232").Value = "Caramelo"
$ws.Range("L232").Value = "Segunda"
$ws.Range("M232").Value = 200
$ws.Range("N232").Value = 18000
$ws.Range("O232").Value = 18000
$ws.Range("P232").Value = 18000
$ws.Range("Q232").Value = "$/caja 14 unidades"
$ws.Range("R232").Value = "Ecuador"
$ws.Range("S232").Value = 1286
$ws.Range("T232").Value = 14

# Make sure the date column keeps the same date number format used by the
# rest of the column (style index 2 in the original workbook).
$ws.Range("D232").NumberFormat = $ws.Range("D233").NumberFormat
